$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.051
$ws.Range("E2").Value = 0.056
$ws.Range("F2").Value = 0.052

$ws.Range("D3").Value = 71.554
$ws.Range("E3").Value = 71.786
$ws.Range("F3").Value = 71.325

$ws.Range("D4").Value = 0.058
$ws.Range("E4").Value = 0.06
$ws.Range("F4").Value = 0.059

$ws.Range("D5").Value = 399.85
$ws.Range("E5").Value = 417.367
$ws.Range("F5").Value = 402.647
$ws.Range("H5").Value = 1169847

$ws.Range("D6").Value = 97.304
$ws.Range("E6").Value = 97.596
$ws.Range("F6").Value = 97.357
$ws.Range("H6").Value = 487802

$ws.Range("D7").Value = 230.154
$ws.Range("E7").Value = 234.9
$ws.Range("F7").Value = 230.282
$ws.Range("H7").Value = 639058

$ws.Range("D8").Value = 843.998
$ws.Range("E8").Value = 860.997
$ws.Range("F8").Value = 848.513
$ws.Range("H8").Value = 1164380

$ws.Range("D9").Value = 569.706
$ws.Range("E9").Value = 575.672
$ws.Range("F9").Value = 570.758
$ws.Range("H9").Value = 1164380

$ws.Range("D10").Value = 565.306
$ws.Range("E10").Value = 570.357
$ws.Range("F10").Value = 564.075
$ws.Range("H10").Value = 1164380

$ws.Range("D11").Value = 815.731
$ws.Range("E11").Value = 833.903
$ws.Range("F11").Value = 820.753
$ws.Range("H11").Value = 1164380

$ws.Range("D12").Value = 485.135
$ws.Range("E12").Value = 493.261
$ws.Range("F12").Value = 486.366
$ws.Range("H12").Value = 1164380

$ws.Range("D13").Value = 481.066
$ws.Range("E13").Value = 507.615
$ws.Range("F13").Value = 486.097
$ws.Range("H13").Value = 1164380

$ws.Range("D14").Value = 411.758
$ws.Range("E14").Value = 417.132
$ws.Range("F14").Value = 412.577
$ws.Range("H14").Value = 1164380

$ws.Range("D15").Value = 1061.979
$ws.Range("E15").Value = 1073.709
$ws.Range("F15").Value = 1062.12
$ws.Range("H15").Value = 1164380

$ws.Range("D16").Value = 969.648
$ws.Range("E16").Value = 984.454
$ws.Range("F16").Value = 974.478
$ws.Range("H16").Value = 1164380

$ws.Range("D17").Value = 633.186
$ws.Range("E17").Value = 653.856
$ws.Range("F17").Value = 634.716
$ws.Range("H17").Value = 1164380

$ws.Range("D18").Value = 813.162
$ws.Range("E18").Value = 820.707
$ws.Range("F18").Value = 813.252
$ws.Range("H18").Value = 1164380

$ws.Range("D19").Value = 583.098
$ws.Range("E19").Value = 587.112
$ws.Range("F19").Value = 582.279
$ws.Range("H19").Value = 1164380

$ws.Range("D20").Value = 573.491
$ws.Range("E20").Value = 579.384
$ws.Range("F20").Value = 574.231
$ws.Range("H20").Value = 1164380

$ws.Range("D21").Value = 856.395
$ws.Range("E21").Value = 872.314
$ws.Range("F21").Value = 858.476
$ws.Range("H21").Value = 1164380

$ws.Range("D22").Value = 0.057
$ws.Range("E22").Value = 0.061
$ws.Range("F22").Value = 0.058

$ws.Range("D23").Value = 51.807
$ws.Range("E23").Value = 53.314
$ws.Range("F23").Value = 51.98
$ws.Range("H23").Value = 267897

$ws.Range("D24").Value = 0.047
$ws.Range("E24").Value = 0.05
$ws.Range("F24").Value = 0.048

$ws.Range("D25").Value = 222.267
$ws.Range("E25").Value = 224.465
$ws.Range("F25").Value = 222.373
$ws.Range("H25").Value = 366037

$ws.Range("D26").Value = 64.31100000000001
$ws.Range("E26").Value = 64.642
$ws.Range("F26").Value = 64.3
$ws.Range("H26").Value = 299068

$ws.Range("D27").Value = 108.478
$ws.Range("E27").Value = 108.956
$ws.Range("F27").Value = 108.29
$ws.Range("H27").Value = 334950

$ws.Range("D28").Value = 224.019
$ws.Range("E28").Value = 228.328
$ws.Range("F28").Value = 223.853
$ws.Range("H28").Value = 362971

$ws.Range("D29").Value = 47.167
$ws.Range("E29").Value = 49.017
$ws.Range("F29").Value = 47.499
$ws.Range("H29").Value = 269087

$ws.Range("D30").Value = 0.011
$ws.Range("E30").Value = 0.012
$ws.Range("F30").Value = 0.011

$ws.Range("D31").Value = 163.378
$ws.Range("E31").Value = 168.793
$ws.Range("F31").Value = 163.771
$ws.Range("H31").Value = 364322

$ws.Range("D32").Value = 160.575
$ws.Range("E32").Value = 168.966
$ws.Range("F32").Value = 162.046
$ws.Range("H32").Value = 361340

$ws.Range("D33").Value = 303.868
$ws.Range("E33").Value = 309.193
$ws.Range("F33").Value = 305.641
$ws.Range("H33").Value = 361340

$ws.Range("D34").Value = 165.032
$ws.Range("E34").Value = 181.473
$ws.Range("F34").Value = 169.158
$ws.Range("H34").Value = 361340

$ws.Range("D35").Value = 654.438
$ws.Range("E35").Value = 671.581
$ws.Range("F35").Value = 653.268
$ws.Range("H35").Value = 361340

$ws.Range("D36").Value = 182.75
$ws.Range("E36").Value = 207.362
$ws.Range("F36").Value = 188.383
$ws.Range("H36").Value = 361340

$ws.Range("D37").Value = 382.212
$ws.Range("E37").Value = 391.605
$ws.Range("F37").Value = 383.677
$ws.Range("H37").Value = 361340

$ws.Range("D38").Value = 389.146
$ws.Range("E38").Value = 401.245
$ws.Range("F38").Value = 384.547
$ws.Range("H38").Value = 361340

$ws.Range("D39").Value = 267.505
$ws.Range("E39").Value = 271.305
$ws.Range("F39").Value = 267.449
$ws.Range("H39").Value = 361340

$ws.Range("D40").Value = 266.553
$ws.Range("E40").Value = 272.727
$ws.Range("F40").Value = 267.631
$ws.Range("H40").Value = 361340

$ws.Range("D41").Value = 322.812
$ws.Range("E41").Value = 329.53
$ws.Range("F41").Value = 322.725
$ws.Range("H41").Value = 361340
